# Updated cryptos list on Tue Oct 22 07:29:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to stay text (matches original inlineStr cells) so that
    # numeric-looking strings (e.g. "596.73", "0.160", "1.00") are not
    # silently converted into numbers / lose formatting by Excel.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "67.638.13"
Set-TextValue $ws.Range("E2") "  -1.51%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.655.88"
Set-TextValue $ws.Range("E3") "  -2.70%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  +0.00%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "596.73"
Set-TextValue $ws.Range("E5") "  -1.95%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "168.25"
Set-TextValue $ws.Range("E6") "  -0.78%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  +0.01%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.547"
Set-TextValue $ws.Range("E8") "  -0.27%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.655.33"
Set-TextValue $ws.Range("E9") "  -2.71%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("E10") "  -1.11%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.160"
Set-TextValue $ws.Range("E11") "  +2.02%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("E12") "  -0.14%  "

# Row 13 - Toncoin
Set-TextValue $ws.Range("D13") "5.28"
Set-TextValue $ws.Range("E13") "  -1.28%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "28.13"
Set-TextValue $ws.Range("E14") "  -2.25%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.137.86"
Set-TextValue $ws.Range("E15") "  -2.79%  "

# Row 16 - ShibaInu
Set-TextValue $ws.Range("E16") "  -3.19%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "67.578.53"
Set-TextValue $ws.Range("E17") "  -1.61%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.642.55"
Set-TextValue $ws.Range("E18") "  -3.62%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "12.17"
Set-TextValue $ws.Range("E19") "  +1.99%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "8.16"
Set-TextValue $ws.Range("E20") "  +6.47%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "363.87"
Set-TextValue $ws.Range("E21") "  -3.29%  "

# Row 22 - Polkadot
Set-TextValue $ws.Range("D22") "4.42"
Set-TextValue $ws.Range("E22") "  -2.17%  "

# Row 23 - NEARProtocol
Set-TextValue $ws.Range("D23") "4.81"
Set-TextValue $ws.Range("E23") "  -3.82%  "

# Row 24 - Aptos
Set-TextValue $ws.Range("D24") "11.05"
Set-TextValue $ws.Range("E24") "  +9.17%  "

# Row 25 - SuiNetwork
Set-TextValue $ws.Range("E25") "  -4.22%  "

# Rows 26 & 27 - Dai and Litecoin swap places, with updated price/volume
Set-TextValue $ws.Range("B26") "Litecoin"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D26") "71.43"
Set-TextValue $ws.Range("E26") "  -3.12%  "

Set-TextValue $ws.Range("B27") "Dai"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D27") "1.00"
Set-TextValue $ws.Range("E27") "  +0.05%  "

# Row 28 - WrappedeETH
Set-TextValue $ws.Range("D28") "2.794.33"
Set-TextValue $ws.Range("E28") "  -2.63%  "

# Row 29 - PEPE
Set-TextValue $ws.Range("E29") "  -2.68%  "

# Row 30 - Binance-PegBSC-USD
Set-TextValue $ws.Range("E30") "  -0.20%  "

# Row 31 - Bittensor
Set-TextValue $ws.Range("D31") "559.38"
Set-TextValue $ws.Range("E31") "  -5.36%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("E32") "  -3.11%  "

# Row 33 - Fetch.AI
Set-TextValue $ws.Range("E33") "  -3.55%  "

# Row 34 - PancakeSwap
Set-TextValue $ws.Range("E34") "  -3.26%  "

# Row 35 - Kaspa
Set-TextValue $ws.Range("D35") "0.135"
Set-TextValue $ws.Range("E35") "  +1.86%  "

# Row 37 - ImmutableX
Set-TextValue $ws.Range("D37") "1.56"
Set-TextValue $ws.Range("E37") "  -4.06%  "

# Row 38 - Monero
Set-TextValue $ws.Range("D38") "158.09"
Set-TextValue $ws.Range("E38") "  -3.23%  "

# Row 39 - EthereumClassic
Set-TextValue $ws.Range("E39") "  -2.88%  "

# Row 40 - PolygonEcosystemToken
Set-TextValue $ws.Range("E40") "  -2.00%  "

# Row 41 - RenderToken
Set-TextValue $ws.Range("D41") "5.33"
Set-TextValue $ws.Range("E41") "  -2.86%  "

# Row 42 - Stacks
Set-TextValue $ws.Range("E42") "  -3.70%  "

# Row 43 - WhiteBITCoin
Set-TextValue $ws.Range("E43") "  -0.22%  "

# Row 44 - dogwifhat
Set-TextValue $ws.Range("E44") "  -4.54%  "

# Row 45 - USDe
Set-TextValue $ws.Range("E45") "  +0.06%  "

# Row 46 - OKB
Set-TextValue $ws.Range("D46") "40.35"
Set-TextValue $ws.Range("E46") "  -1.56%  "

# Row 47 - BabyDogeCoin
Set-TextValue $ws.Range("E47") "  -3.86%  "

# Row 48 - ARBITRUM
Set-TextValue $ws.Range("D48") "0.597"
Set-TextValue $ws.Range("E48") "  -1.22%  "

# Row 49 - Aave
Set-TextValue $ws.Range("D49") "154.68"
Set-TextValue $ws.Range("E49") "  -0.67%  "

# Row 50 - Filecoin
Set-TextValue $ws.Range("E50") "  -1.24%  "

# Row 51 - Optimism
Set-TextValue $ws.Range("E51") "  -2.91%  "
